$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.018.06'
$ws.Range('D3').Value = '2.339.52'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.59%  '
$ws.Range('E7').Value = '  -3.66%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.512'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.33'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0801'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.85'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.88'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.73%  '
$ws.Range('D16').Value = '2.316.57'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = '42.932.44'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.84'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.00%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0913'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.17'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.42%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.49%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.43%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '163.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.95%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.89%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.68'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.45%  '
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.61'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.67%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0729'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.86'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.88%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.102'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.113'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.53%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.007.06'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0286'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.88'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.17'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.92%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '56.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.562.79'
$ws.Range('E51').Value = '  +1.12%  '
